$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-08 Wednesday" "2025-01-14 Tuesday"

Replace-Text "97÷6=" "63÷3="
Replace-Text "10÷2=" "30÷6="
Replace-Text "38÷9=" "89÷5="
Replace-Text "82÷8=" "72÷4="
Replace-Text "68÷5=" "72÷6="
Replace-Text "88÷5=" "72÷7="
Replace-Text "86÷9=" "22÷6="
Replace-Text "78÷5=" "75÷5="
Replace-Text "88÷6=" "39÷6="
Replace-Text "62÷5=" "70÷3="
Replace-Text "14÷8=" "24÷2="
Replace-Text "16÷4=" "97÷4="
Replace-Text "55÷7=" "55÷8="
Replace-Text "90÷6=" "65÷7="
Replace-Text "47÷7=" "71÷5="
Replace-Text "15÷2=" "88÷9="
Replace-Text "64÷8=" "81÷3="
Replace-Text "46÷6=" "69÷8="
Replace-Text "80÷9=" "12÷8="
Replace-Text "49÷9=" "74÷8="
Replace-Text "82÷3=" "75÷9="
Replace-Text "36÷7=" "43÷7="
Replace-Text "65÷4=" "21÷4="
Replace-Text "89÷7=" "88÷7="
Replace-Text "86÷5=" "90÷6="
